# Automatic update of files.
# Applies the per-row edits described by the commit diff:
#  - Taxonsorteringsordning (col B) bumped 79243 -> 79244 on affected rows
#  - Rows 20/21 swap their full record content (two observations re-ordered)
#  - Rows 24/25 swap their full record content (two observations re-ordered)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple rows: bump Taxonsorteringsordning from 79243 to 79244 ---
$simpleRows = @(3, 6, 7, 8, 10, 11, 12, 13, 16, 17, 18, 22, 23, 26)
foreach ($r in $simpleRows) {
    $ws.Range("B$r").Value = 79244
}

# --- Row 20: rebuild full record (swapped with its counterpart row) ---
foreach ($addr in @("L20", "M20")) { $ws.Range($addr).ClearContents() | Out-Null }
$ws.Range("A20").Value = 131022915
$ws.Range("B20").Value = 79244
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("I20").Value = "x"
$ws.Range("I20").Value = ""
$ws.Range("J20").Value = "x"
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = "x"
$ws.Range("K20").Value = ""
$ws.Range("N20").Value = "x"
$ws.Range("N20").Value = ""
$ws.Range("P20").Value = "Trossbygget, Dlr"
$ws.Range("Q20").Value = 477093
$ws.Range("R20").Value = 6788924
$ws.Range("S20").Value = 10
$ws.Range("T20").Value = "Dalarna"
$ws.Range("U20").Value = "Orsa"
$ws.Range("V20").Value = "Dalarna"
$ws.Range("W20").Value = "Orsa"
$ws.Range("AD20").Value = $false
$ws.Range("AE20").Value = $false
$ws.Range("AF20").Value = "x"
$ws.Range("AF20").Value = ""
$ws.Range("AG20").Value = $false
$ws.Range("AT20").Value = "x"
$ws.Range("AT20").Value = ""
$ws.Range("AW20").Value = "Håkan Thenander"
$ws.Range("AX20").Value = "Håkan Thenander"
$ws.Range("AY20").Value = "x"
$ws.Range("AY20").Value = ""

# --- Row 21: rebuild full record (swapped with its counterpart row) ---
foreach ($addr in @("J21", "AF21")) { $ws.Range($addr).ClearContents() | Out-Null }
$ws.Range("A21").Value = 131022892
$ws.Range("B21").Value = 57884
$ws.Range("D21").Value = "NT"
$ws.Range("E21").Value = 100109
$ws.Range("F21").Value = "Tretåig hackspett"
$ws.Range("G21").Value = "Picoides tridactylus"
$ws.Range("H21").Value = "(Linnaeus, 1758)"
$ws.Range("I21").Value = "x"
$ws.Range("I21").Value = ""
$ws.Range("K21").Value = "x"
$ws.Range("K21").Value = ""
$ws.Range("L21").Value = "x"
$ws.Range("L21").Value = ""
$ws.Range("M21").Value = "färska spår"
$ws.Range("N21").Value = "x"
$ws.Range("N21").Value = ""
$ws.Range("P21").Value = "Trossbygget, Dlr"
$ws.Range("Q21").Value = 477106
$ws.Range("R21").Value = 6788935
$ws.Range("S21").Value = 10
$ws.Range("T21").Value = "Dalarna"
$ws.Range("U21").Value = "Orsa"
$ws.Range("V21").Value = "Dalarna"
$ws.Range("W21").Value = "Orsa"
$ws.Range("AD21").Value = $false
$ws.Range("AE21").Value = $false
$ws.Range("AG21").Value = $false
$ws.Range("AT21").Value = "x"
$ws.Range("AT21").Value = ""
$ws.Range("AW21").Value = "Håkan Thenander"
$ws.Range("AX21").Value = "Håkan Thenander"
$ws.Range("AY21").Value = "x"
$ws.Range("AY21").Value = ""

# --- Row 24: rebuild full record (swapped with its counterpart row) ---
foreach ($addr in @("Z24", "AB24", "AC24")) { $ws.Range($addr).ClearContents() | Out-Null }
$ws.Range("A24").Value = 131022847
$ws.Range("B24").Value = 57884
$ws.Range("D24").Value = "NT"
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("I24").Value = "x"
$ws.Range("I24").Value = ""
$ws.Range("K24").Value = "x"
$ws.Range("K24").Value = ""
$ws.Range("L24").Value = "x"
$ws.Range("L24").Value = ""
$ws.Range("M24").Value = "färska spår"
$ws.Range("N24").Value = "x"
$ws.Range("N24").Value = ""
$ws.Range("P24").Value = "Trossbygget, Dlr"
$ws.Range("Q24").Value = 477122
$ws.Range("R24").Value = 6788910
$ws.Range("S24").Value = 10
$ws.Range("T24").Value = "Dalarna"
$ws.Range("U24").Value = "Orsa"
$ws.Range("V24").Value = "Dalarna"
$ws.Range("W24").Value = "Orsa"
$ws.Range("AD24").Value = $false
$ws.Range("AE24").Value = $false
$ws.Range("AG24").Value = $false
$ws.Range("AT24").Value = "x"
$ws.Range("AT24").Value = ""
$ws.Range("AW24").Value = "Håkan Thenander"
$ws.Range("AX24").Value = "Håkan Thenander"
$ws.Range("AY24").Value = "x"
$ws.Range("AY24").Value = ""

# --- Row 25: rebuild full record (swapped with its counterpart row) ---
foreach ($addr in @("K25", "L25", "M25", "N25")) { $ws.Range($addr).ClearContents() | Out-Null }
$ws.Range("A25").Value = 131017116
$ws.Range("B25").Value = 79244
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 6425
$ws.Range("F25").Value = "Garnlav"
$ws.Range("G25").Value = "Alectoria sarmentosa"
$ws.Range("H25").Value = "(Ach.) Ach."
$ws.Range("I25").Value = "x"
$ws.Range("I25").Value = ""
$ws.Range("P25").Value = "Trossbygget, Dlr"
$ws.Range("Q25").Value = 477185
$ws.Range("R25").Value = 6789174
$ws.Range("S25").Value = 10
$ws.Range("T25").Value = "Dalarna"
$ws.Range("U25").Value = "Orsa"
$ws.Range("V25").Value = "Dalarna"
$ws.Range("W25").Value = "Orsa"
$ws.Range("Z25").Value = "12:48"
$ws.Range("AB25").Value = "12:48"
$ws.Range("AC25").Value = "Rikligt i området"
$ws.Range("AD25").Value = $false
$ws.Range("AE25").Value = $false
$ws.Range("AG25").Value = $false
$ws.Range("AT25").Value = "x"
$ws.Range("AT25").Value = ""
$ws.Range("AW25").Value = "Håkan Thenander"
$ws.Range("AX25").Value = "Håkan Thenander"
$ws.Range("AY25").Value = "x"
$ws.Range("AY25").Value = ""

Write-Output "done"
